# base_dados_alunos.xlsx
# - Telas editar matricula 1 e 2 adicionadas (novos registros de alunos nas linhas 22-24).
# - Funcionalidade de excluir matricula: a celula AW21 ficava vazia (placeholder) e e removida.
# - Funcionalidade de editar matricula em desenvolvimento (linha 24 ainda tem campos AX/AY/AZ em branco).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Linha 21: remove a celula vazia "AW21" (matricula de registro civil nao informada) ---
$ws.Range("AW21").Value = ""

# --- Linha 22: nova matricula ---
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = '2'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'Renan Rodrigues'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = '3213'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = 'Branca'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = 'Masculino'
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = 'aa'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = 'a'
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = 'aa'
$ws.Range("I22").NumberFormat = "@"
$ws.Range("I22").Value = 'a'
$ws.Range("J22").NumberFormat = "@"
$ws.Range("J22").Value = 'a'
$ws.Range("K22").NumberFormat = "@"
$ws.Range("K22").Value = '2000-01-01'
$ws.Range("L22").NumberFormat = "@"
$ws.Range("L22").Value = 'a'
$ws.Range("M22").NumberFormat = "@"
$ws.Range("M22").Value = 'aa'
$ws.Range("N22").NumberFormat = "@"
$ws.Range("N22").Value = '1'
$ws.Range("O22").NumberFormat = "@"
$ws.Range("O22").Value = 'NÃO'
$ws.Range("P22").NumberFormat = "@"
$ws.Range("P22").Value = 'NÃO'
$ws.Range("Q22").NumberFormat = "@"
$ws.Range("Q22").Value = 'NÃO'
$ws.Range("R22").NumberFormat = "@"
$ws.Range("R22").Value = 'NÃO'
$ws.Range("S22").NumberFormat = "@"
$ws.Range("S22").Value = 'SIM'
$ws.Range("T22").NumberFormat = "@"
$ws.Range("T22").Value = 'NÃO'
$ws.Range("U22").NumberFormat = "@"
$ws.Range("U22").Value = 'NÃO'
$ws.Range("V22").NumberFormat = "@"
$ws.Range("V22").Value = 'NÃO'
$ws.Range("W22").NumberFormat = "@"
$ws.Range("W22").Value = 'NÃO'
$ws.Range("X22").NumberFormat = "@"
$ws.Range("X22").Value = 'NÃO'
$ws.Range("Y22").NumberFormat = "@"
$ws.Range("Y22").Value = 'NÃO'
$ws.Range("Z22").NumberFormat = "@"
$ws.Range("Z22").Value = 'NÃO'
$ws.Range("AA22").NumberFormat = "@"
$ws.Range("AA22").Value = 'NÃO'
$ws.Range("AB22").NumberFormat = "@"
$ws.Range("AB22").Value = 'NÃO'
$ws.Range("AC22").NumberFormat = "@"
$ws.Range("AC22").Value = 'a'
$ws.Range("AD22").NumberFormat = "@"
$ws.Range("AD22").Value = 'a'
$ws.Range("AE22").NumberFormat = "@"
$ws.Range("AE22").Value = '1'
$ws.Range("AF22").NumberFormat = "@"
$ws.Range("AF22").Value = 'a'
$ws.Range("AG22").NumberFormat = "@"
$ws.Range("AG22").Value = '1'
$ws.Range("AH22").NumberFormat = "@"
$ws.Range("AH22").Value = 'Urbana'
$ws.Range("AI22").NumberFormat = "@"
$ws.Range("AI22").Value = '1'
$ws.Range("AJ22").NumberFormat = "@"
$ws.Range("AJ22").Value = 'b@b.b'
$ws.Range("AK22").NumberFormat = "@"
$ws.Range("AK22").Value = 'pai'
$ws.Range("AL22").NumberFormat = "@"
$ws.Range("AL22").Value = 'mae'
$ws.Range("AM22").NumberFormat = "@"
$ws.Range("AM22").Value = 'escola'
$ws.Range("AN22").NumberFormat = "@"
$ws.Range("AN22").Value = '1'
$ws.Range("AO22").NumberFormat = "@"
$ws.Range("AO22").Value = '1/1/2000'
$ws.Range("AP22").NumberFormat = "@"
$ws.Range("AP22").Value = '1'
$ws.Range("AQ22").NumberFormat = "@"
$ws.Range("AQ22").Value = '1/1/2000'
$ws.Range("AR22").NumberFormat = "@"
$ws.Range("AR22").Value = 'Manhã'
$ws.Range("AS22").NumberFormat = "@"
$ws.Range("AS22").Value = '01. Berçário I'
$ws.Range("AT22").NumberFormat = "@"
$ws.Range("AT22").Value = '01 - Do Lar'
$ws.Range("AU22").NumberFormat = "@"
$ws.Range("AU22").Value = 'NÃO'
$ws.Range("AV22").NumberFormat = "@"
$ws.Range("AV22").Value = 'NÃO'
$ws.Range("AW22").NumberFormat = "@"
$ws.Range("AW22").Value = '1'
$ws.Range("BA22").NumberFormat = "@"
$ws.Range("BA22").Value = '2000-01-01'

# --- Linha 23: nova matricula ---
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = '220713'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'nome aluno'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = '123'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = 'Preta'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = 'Feminino'
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = 'UE'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = 'municipio endereco'
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = 'UC'
$ws.Range("I23").NumberFormat = "@"
$ws.Range("I23").Value = 'nome cartorio'
$ws.Range("J23").NumberFormat = "@"
$ws.Range("J23").Value = 'municipio cartorio'
$ws.Range("K23").NumberFormat = "@"
$ws.Range("K23").Value = '2000-01-01'
$ws.Range("L23").NumberFormat = "@"
$ws.Range("L23").Value = 'OE'
$ws.Range("M23").NumberFormat = "@"
$ws.Range("M23").Value = 'UF'
$ws.Range("N23").NumberFormat = "@"
$ws.Range("N23").Value = '758'
$ws.Range("O23").NumberFormat = "@"
$ws.Range("O23").Value = 'NÃO'
$ws.Range("P23").NumberFormat = "@"
$ws.Range("P23").Value = 'NÃO'
$ws.Range("Q23").NumberFormat = "@"
$ws.Range("Q23").Value = 'NÃO'
$ws.Range("R23").NumberFormat = "@"
$ws.Range("R23").Value = 'NÃO'
$ws.Range("S23").NumberFormat = "@"
$ws.Range("S23").Value = 'NÃO'
$ws.Range("T23").NumberFormat = "@"
$ws.Range("T23").Value = 'NÃO'
$ws.Range("U23").NumberFormat = "@"
$ws.Range("U23").Value = 'NÃO'
$ws.Range("V23").NumberFormat = "@"
$ws.Range("V23").Value = 'NÃO'
$ws.Range("W23").NumberFormat = "@"
$ws.Range("W23").Value = 'NÃO'
$ws.Range("X23").NumberFormat = "@"
$ws.Range("X23").Value = 'NÃO'
$ws.Range("Y23").NumberFormat = "@"
$ws.Range("Y23").Value = 'NÃO'
$ws.Range("Z23").NumberFormat = "@"
$ws.Range("Z23").Value = 'NÃO'
$ws.Range("AA23").NumberFormat = "@"
$ws.Range("AA23").Value = 'NÃO'
$ws.Range("AB23").NumberFormat = "@"
$ws.Range("AB23").Value = 'SIM'
$ws.Range("AC23").NumberFormat = "@"
$ws.Range("AC23").Value = 'endereco'
$ws.Range("AD23").NumberFormat = "@"
$ws.Range("AD23").Value = 'complement'
$ws.Range("AE23").NumberFormat = "@"
$ws.Range("AE23").Value = '7'
$ws.Range("AF23").NumberFormat = "@"
$ws.Range("AF23").Value = 'bairro'
$ws.Range("AG23").NumberFormat = "@"
$ws.Range("AG23").Value = '666'
$ws.Range("AH23").NumberFormat = "@"
$ws.Range("AH23").Value = 'Urbana'
$ws.Range("AI23").NumberFormat = "@"
$ws.Range("AI23").Value = '555'
$ws.Range("AJ23").NumberFormat = "@"
$ws.Range("AJ23").Value = 'email@email.com'
$ws.Range("AK23").NumberFormat = "@"
$ws.Range("AK23").Value = 'Pai'
$ws.Range("AL23").NumberFormat = "@"
$ws.Range("AL23").Value = 'Mãe'
$ws.Range("AM23").NumberFormat = "@"
$ws.Range("AM23").Value = 'escola'
$ws.Range("AN23").NumberFormat = "@"
$ws.Range("AN23").Value = '5'
$ws.Range("AO23").NumberFormat = "@"
$ws.Range("AO23").Value = '1/1/2000'
$ws.Range("AP23").NumberFormat = "@"
$ws.Range("AP23").Value = '9'
$ws.Range("AQ23").NumberFormat = "@"
$ws.Range("AQ23").Value = '1/1/2000'
$ws.Range("AR23").NumberFormat = "@"
$ws.Range("AR23").Value = 'Intermediário'
$ws.Range("AS23").NumberFormat = "@"
$ws.Range("AS23").Value = '02. Berçário II'
$ws.Range("AT23").NumberFormat = "@"
$ws.Range("AT23").Value = '04 - Escola Particular'
$ws.Range("AU23").NumberFormat = "@"
$ws.Range("AU23").Value = 'SIM'
$ws.Range("AV23").NumberFormat = "@"
$ws.Range("AV23").Value = 'NÃO'
$ws.Range("AW23").NumberFormat = "@"
$ws.Range("AW23").Value = '321'
$ws.Range("BA23").NumberFormat = "@"
$ws.Range("BA23").Value = '2000-01-01'

# --- Linha 24: nova matricula (edicao em desenvolvimento) ---
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = '1'
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = 'nome'
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = '1'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = 'Branca'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = 'Masculino'
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = 'aa'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = 'a'
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = 'aa'
$ws.Range("I24").NumberFormat = "@"
$ws.Range("I24").Value = 'a'
$ws.Range("J24").NumberFormat = "@"
$ws.Range("J24").Value = 'a'
$ws.Range("K24").NumberFormat = "@"
$ws.Range("K24").Value = '2000-01-01'
$ws.Range("L24").NumberFormat = "@"
$ws.Range("L24").Value = 'aa'
$ws.Range("M24").NumberFormat = "@"
$ws.Range("M24").Value = 'aa'
$ws.Range("N24").NumberFormat = "@"
$ws.Range("N24").Value = '1'
$ws.Range("O24").NumberFormat = "@"
$ws.Range("O24").Value = 'SIM'
$ws.Range("P24").NumberFormat = "@"
$ws.Range("P24").Value = 'SIM'
$ws.Range("Q24").NumberFormat = "@"
$ws.Range("Q24").Value = 'SIM'
$ws.Range("R24").NumberFormat = "@"
$ws.Range("R24").Value = 'SIM'
$ws.Range("S24").NumberFormat = "@"
$ws.Range("S24").Value = 'SIM'
$ws.Range("T24").NumberFormat = "@"
$ws.Range("T24").Value = 'SIM'
$ws.Range("U24").NumberFormat = "@"
$ws.Range("U24").Value = 'SIM'
$ws.Range("V24").NumberFormat = "@"
$ws.Range("V24").Value = 'SIM'
$ws.Range("W24").NumberFormat = "@"
$ws.Range("W24").Value = 'SIM'
$ws.Range("X24").NumberFormat = "@"
$ws.Range("X24").Value = 'SIM'
$ws.Range("Y24").NumberFormat = "@"
$ws.Range("Y24").Value = 'SIM'
$ws.Range("Z24").NumberFormat = "@"
$ws.Range("Z24").Value = 'SIM'
$ws.Range("AA24").NumberFormat = "@"
$ws.Range("AA24").Value = 'SIM'
$ws.Range("AB24").NumberFormat = "@"
$ws.Range("AB24").Value = 'SIM'
$ws.Range("AC24").NumberFormat = "@"
$ws.Range("AC24").Value = 'a'
$ws.Range("AD24").NumberFormat = "@"
$ws.Range("AD24").Value = 'a'
$ws.Range("AE24").NumberFormat = "@"
$ws.Range("AE24").Value = '1'
$ws.Range("AF24").NumberFormat = "@"
$ws.Range("AF24").Value = 'a'
$ws.Range("AG24").NumberFormat = "@"
$ws.Range("AG24").Value = '1'
$ws.Range("AH24").NumberFormat = "@"
$ws.Range("AH24").Value = 'Rural'
$ws.Range("AI24").NumberFormat = "@"
$ws.Range("AI24").Value = '1'
$ws.Range("AJ24").NumberFormat = "@"
$ws.Range("AJ24").Value = 'a'
$ws.Range("AK24").NumberFormat = "@"
$ws.Range("AK24").Value = 'pai'
$ws.Range("AL24").NumberFormat = "@"
$ws.Range("AL24").Value = 'mae'
$ws.Range("AM24").NumberFormat = "@"
$ws.Range("AM24").Value = 'a'
$ws.Range("AN24").NumberFormat = "@"
$ws.Range("AN24").Value = '1'
$ws.Range("AO24").NumberFormat = "@"
$ws.Range("AO24").Value = '1/1/2000'
$ws.Range("AP24").NumberFormat = "@"
$ws.Range("AP24").Value = '1'
$ws.Range("AQ24").NumberFormat = "@"
$ws.Range("AQ24").Value = '1/1/2000'
$ws.Range("AR24").NumberFormat = "@"
$ws.Range("AR24").Value = 'Manhã'
$ws.Range("AS24").NumberFormat = "@"
$ws.Range("AS24").Value = '01. Berçário I'
$ws.Range("AT24").NumberFormat = "@"
$ws.Range("AT24").Value = '02 - Escola Municipal'
$ws.Range("AU24").NumberFormat = "@"
$ws.Range("AU24").Value = 'NÃO'
$ws.Range("AV24").NumberFormat = "@"
$ws.Range("AV24").Value = 'NÃO'
$ws.Range("AW24").NumberFormat = "@"
$ws.Range("AW24").Value = '1'
$ws.Range("AX24").NumberFormat = "@"
$ws.Range("AX24").Value = "'"
$ws.Range("AY24").NumberFormat = "@"
$ws.Range("AY24").Value = "'"
$ws.Range("AZ24").NumberFormat = "@"
$ws.Range("AZ24").Value = "'"
$ws.Range("BA24").NumberFormat = "@"
$ws.Range("BA24").Value = '2000-01-01'

